$p = $ppt.ActivePresentation

# --- Slide 3 ("Aspectos Economicos / Quanto custa o problema?") ---
# CaixaDeTexto 2 (shape 5): shrink the box and replace the last bullet's
# long sentence with the short heading "Fitoterapicos e custos".
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(5)
$tr3 = $sh3.TextFrame.TextRange
$tr3.Paragraphs(5,1).Text = "Fitoterápicos e custos  "
$sh3.Height = 225.3797

# --- Slide 6 ("Demanda de Mercado") ---
# CaixaDeTexto 2 (shape 5): shrink the box, shorten the first and third
# bullets to short headings, and drop one trailing empty paragraph.
$s6 = $p.Slides.Item(6)
$sh6 = $s6.Shapes.Item(5)
$tr6 = $sh6.TextFrame.TextRange
$tr6.Paragraphs(1,1).Text = "o mercado mundial"
$tr6.Paragraphs(5,1).Text = "O mercado brasileiro"
$tr6.Paragraphs(7,1).Delete()
$sh6.Height = 159.9469
